# Auto-generated Excel COM-interop script
# Applies the refreshed crypto market-data snapshot: updated Price /
# Volume(1h) / Hora columns for every listed coin, plus the rotation of
# coin name+link pairs across rows 6-18 that came from the upstream list
# being re-sorted before this run's scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "333.87"
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "7.52%"
$ws.Cells.Item(2, 5).ClearFormats()
$ws.Cells.Item(2, 7).NumberFormat = "@"
$ws.Cells.Item(2, 7).Value = "6"
$ws.Cells.Item(2, 7).ClearFormats()

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "40.63"
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "7.82%"
$ws.Cells.Item(3, 5).ClearFormats()
$ws.Cells.Item(3, 7).NumberFormat = "@"
$ws.Cells.Item(3, 7).Value = "6"
$ws.Cells.Item(3, 7).ClearFormats()

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "5.272"
$ws.Cells.Item(4, 4).ClearFormats()
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "1.98%"
$ws.Cells.Item(4, 5).ClearFormats()
$ws.Cells.Item(4, 7).NumberFormat = "@"
$ws.Cells.Item(4, 7).Value = "6"
$ws.Cells.Item(4, 7).ClearFormats()

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "0.08095"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "2.40%"
$ws.Cells.Item(5, 5).ClearFormats()
$ws.Cells.Item(5, 7).NumberFormat = "@"
$ws.Cells.Item(5, 7).Value = "6"
$ws.Cells.Item(5, 7).ClearFormats()

# Row 6
$ws.Cells.Item(6, 2).Value = "KuCoinToken"
$ws.Cells.Item(6, 3).Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "8.672"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "4.88%"
$ws.Cells.Item(6, 5).ClearFormats()
$ws.Cells.Item(6, 7).NumberFormat = "@"
$ws.Cells.Item(6, 7).Value = "6"
$ws.Cells.Item(6, 7).ClearFormats()

# Row 7
$ws.Cells.Item(7, 2).Value = "FTXToken"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "1.913"
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = "0.18%"
$ws.Cells.Item(7, 5).ClearFormats()
$ws.Cells.Item(7, 7).NumberFormat = "@"
$ws.Cells.Item(7, 7).Value = "6"
$ws.Cells.Item(7, 7).ClearFormats()

# Row 8
$ws.Cells.Item(8, 2).Value = "BTSEToken"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "2.957"
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = "-1.40%"
$ws.Cells.Item(8, 5).ClearFormats()
$ws.Cells.Item(8, 7).NumberFormat = "@"
$ws.Cells.Item(8, 7).Value = "6"
$ws.Cells.Item(8, 7).ClearFormats()

# Row 9
$ws.Cells.Item(9, 2).Value = "MXToken"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.9366"
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = "-0.17%"
$ws.Cells.Item(9, 5).ClearFormats()
$ws.Cells.Item(9, 7).NumberFormat = "@"
$ws.Cells.Item(9, 7).Value = "6"
$ws.Cells.Item(9, 7).ClearFormats()

# Row 10
$ws.Cells.Item(10, 2).Value = "LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.1353"
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "21.55%"
$ws.Cells.Item(10, 5).ClearFormats()
$ws.Cells.Item(10, 7).NumberFormat = "@"
$ws.Cells.Item(10, 7).Value = "6"
$ws.Cells.Item(10, 7).ClearFormats()

# Row 11
$ws.Cells.Item(11, 2).Value = "WazirX"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.1970"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = "0.97%"
$ws.Cells.Item(11, 5).ClearFormats()
$ws.Cells.Item(11, 7).NumberFormat = "@"
$ws.Cells.Item(11, 7).Value = "6"
$ws.Cells.Item(11, 7).ClearFormats()

# Row 12
$ws.Cells.Item(12, 2).Value = "MandalaExchangeToken"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.09130"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = "0.35%"
$ws.Cells.Item(12, 5).ClearFormats()
$ws.Cells.Item(12, 7).NumberFormat = "@"
$ws.Cells.Item(12, 7).Value = "6"
$ws.Cells.Item(12, 7).ClearFormats()

# Row 13
$ws.Cells.Item(13, 2).Value = "BitrueCoin"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.03434"
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "3.05%"
$ws.Cells.Item(13, 5).ClearFormats()
$ws.Cells.Item(13, 7).NumberFormat = "@"
$ws.Cells.Item(13, 7).Value = "6"
$ws.Cells.Item(13, 7).ClearFormats()

# Row 14
$ws.Cells.Item(14, 2).Value = "BitMartToken"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.09569"
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = "-0.42%"
$ws.Cells.Item(14, 5).ClearFormats()
$ws.Cells.Item(14, 7).NumberFormat = "@"
$ws.Cells.Item(14, 7).Value = "6"
$ws.Cells.Item(14, 7).ClearFormats()

# Row 15
$ws.Cells.Item(15, 2).Value = "BitForexToken"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.001400"
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = "0.51%"
$ws.Cells.Item(15, 5).ClearFormats()
$ws.Cells.Item(15, 7).NumberFormat = "@"
$ws.Cells.Item(15, 7).Value = "6"
$ws.Cells.Item(15, 7).ClearFormats()

# Row 16
$ws.Cells.Item(16, 2).Value = "TigerCash"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.006560"
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = "14.50%"
$ws.Cells.Item(16, 5).ClearFormats()
$ws.Cells.Item(16, 7).NumberFormat = "@"
$ws.Cells.Item(16, 7).Value = "6"
$ws.Cells.Item(16, 7).ClearFormats()

# Row 17
$ws.Cells.Item(17, 2).Value = "LEO"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "3.359"
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = "-6.56%"
$ws.Cells.Item(17, 5).ClearFormats()
$ws.Cells.Item(17, 7).NumberFormat = "@"
$ws.Cells.Item(17, 7).Value = "6"
$ws.Cells.Item(17, 7).ClearFormats()

# Row 18
$ws.Cells.Item(18, 2).Value = "GateToken"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "4.532"
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = "2.29%"
$ws.Cells.Item(18, 5).ClearFormats()
$ws.Cells.Item(18, 7).NumberFormat = "@"
$ws.Cells.Item(18, 7).Value = "6"
$ws.Cells.Item(18, 7).ClearFormats()

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.3524"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = "3.27%"
$ws.Cells.Item(19, 5).ClearFormats()
$ws.Cells.Item(19, 7).NumberFormat = "@"
$ws.Cells.Item(19, 7).Value = "6"
$ws.Cells.Item(19, 7).ClearFormats()

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "6.476"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = "0.84%"
$ws.Cells.Item(20, 5).ClearFormats()
$ws.Cells.Item(20, 7).NumberFormat = "@"
$ws.Cells.Item(20, 7).Value = "6"
$ws.Cells.Item(20, 7).ClearFormats()

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.1324"
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = "3.31%"
$ws.Cells.Item(21, 5).ClearFormats()
$ws.Cells.Item(21, 7).NumberFormat = "@"
$ws.Cells.Item(21, 7).Value = "6"
$ws.Cells.Item(21, 7).ClearFormats()

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.2572"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = "2.01%"
$ws.Cells.Item(22, 5).ClearFormats()
$ws.Cells.Item(22, 7).NumberFormat = "@"
$ws.Cells.Item(22, 7).Value = "6"
$ws.Cells.Item(22, 7).ClearFormats()

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.04442"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = "1.00%"
$ws.Cells.Item(23, 5).ClearFormats()
$ws.Cells.Item(23, 7).NumberFormat = "@"
$ws.Cells.Item(23, 7).Value = "6"
$ws.Cells.Item(23, 7).ClearFormats()

# Row 24
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = "-0.89%"
$ws.Cells.Item(24, 5).ClearFormats()
$ws.Cells.Item(24, 7).NumberFormat = "@"
$ws.Cells.Item(24, 7).Value = "6"
$ws.Cells.Item(24, 7).ClearFormats()

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.004326"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = "-5.67%"
$ws.Cells.Item(25, 5).ClearFormats()
$ws.Cells.Item(25, 7).NumberFormat = "@"
$ws.Cells.Item(25, 7).Value = "6"
$ws.Cells.Item(25, 7).ClearFormats()

# Row 26
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = "-5.19%"
$ws.Cells.Item(26, 5).ClearFormats()
$ws.Cells.Item(26, 7).NumberFormat = "@"
$ws.Cells.Item(26, 7).Value = "6"
$ws.Cells.Item(26, 7).ClearFormats()

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.0003996"
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = "0.11%"
$ws.Cells.Item(27, 5).ClearFormats()
$ws.Cells.Item(27, 7).NumberFormat = "@"
$ws.Cells.Item(27, 7).Value = "6"
$ws.Cells.Item(27, 7).ClearFormats()

# Row 28
$ws.Cells.Item(28, 7).NumberFormat = "@"
$ws.Cells.Item(28, 7).Value = "6"
$ws.Cells.Item(28, 7).ClearFormats()

# Row 29
$ws.Cells.Item(29, 7).NumberFormat = "@"
$ws.Cells.Item(29, 7).Value = "6"
$ws.Cells.Item(29, 7).ClearFormats()

# Row 30
$ws.Cells.Item(30, 7).NumberFormat = "@"
$ws.Cells.Item(30, 7).Value = "6"
$ws.Cells.Item(30, 7).ClearFormats()

# Row 31
$ws.Cells.Item(31, 7).NumberFormat = "@"
$ws.Cells.Item(31, 7).Value = "6"
$ws.Cells.Item(31, 7).ClearFormats()

# Row 32
$ws.Cells.Item(32, 7).NumberFormat = "@"
$ws.Cells.Item(32, 7).Value = "6"
$ws.Cells.Item(32, 7).ClearFormats()

# Row 33
$ws.Cells.Item(33, 7).NumberFormat = "@"
$ws.Cells.Item(33, 7).Value = "6"
$ws.Cells.Item(33, 7).ClearFormats()

# Row 34
$ws.Cells.Item(34, 7).NumberFormat = "@"
$ws.Cells.Item(34, 7).Value = "6"
$ws.Cells.Item(34, 7).ClearFormats()

# Row 35
$ws.Cells.Item(35, 7).NumberFormat = "@"
$ws.Cells.Item(35, 7).Value = "6"
$ws.Cells.Item(35, 7).ClearFormats()

# Row 36
$ws.Cells.Item(36, 7).NumberFormat = "@"
$ws.Cells.Item(36, 7).Value = "6"
$ws.Cells.Item(36, 7).ClearFormats()

# Row 37
$ws.Cells.Item(37, 7).NumberFormat = "@"
$ws.Cells.Item(37, 7).Value = "6"
$ws.Cells.Item(37, 7).ClearFormats()

# Row 38
$ws.Cells.Item(38, 7).NumberFormat = "@"
$ws.Cells.Item(38, 7).Value = "6"
$ws.Cells.Item(38, 7).ClearFormats()

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.02509"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = "12.21%"
$ws.Cells.Item(39, 5).ClearFormats()
$ws.Cells.Item(39, 7).NumberFormat = "@"
$ws.Cells.Item(39, 7).Value = "6"
$ws.Cells.Item(39, 7).ClearFormats()

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.05220"
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = "1.79%"
$ws.Cells.Item(40, 5).ClearFormats()
$ws.Cells.Item(40, 7).NumberFormat = "@"
$ws.Cells.Item(40, 7).Value = "6"
$ws.Cells.Item(40, 7).ClearFormats()

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.007663"
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = "2.69%"
$ws.Cells.Item(41, 5).ClearFormats()
$ws.Cells.Item(41, 7).NumberFormat = "@"
$ws.Cells.Item(41, 7).Value = "6"
$ws.Cells.Item(41, 7).ClearFormats()

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.1431"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = "5.77%"
$ws.Cells.Item(42, 5).ClearFormats()
$ws.Cells.Item(42, 7).NumberFormat = "@"
$ws.Cells.Item(42, 7).Value = "6"
$ws.Cells.Item(42, 7).ClearFormats()

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.009046"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = "3.35%"
$ws.Cells.Item(43, 5).ClearFormats()
$ws.Cells.Item(43, 7).NumberFormat = "@"
$ws.Cells.Item(43, 7).Value = "6"
$ws.Cells.Item(43, 7).ClearFormats()

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.002170"
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = "1.83%"
$ws.Cells.Item(44, 5).ClearFormats()
$ws.Cells.Item(44, 7).NumberFormat = "@"
$ws.Cells.Item(44, 7).Value = "6"
$ws.Cells.Item(44, 7).ClearFormats()

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.008989"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = "-3.54%"
$ws.Cells.Item(45, 5).ClearFormats()
$ws.Cells.Item(45, 7).NumberFormat = "@"
$ws.Cells.Item(45, 7).Value = "6"
$ws.Cells.Item(45, 7).ClearFormats()

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.00006631"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = "0.14%"
$ws.Cells.Item(46, 5).ClearFormats()
$ws.Cells.Item(46, 7).NumberFormat = "@"
$ws.Cells.Item(46, 7).Value = "6"
$ws.Cells.Item(46, 7).ClearFormats()

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.00000000751"
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = "0.09%"
$ws.Cells.Item(47, 5).ClearFormats()
$ws.Cells.Item(47, 7).NumberFormat = "@"
$ws.Cells.Item(47, 7).Value = "6"
$ws.Cells.Item(47, 7).ClearFormats()

# Row 48
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = "16.83%"
$ws.Cells.Item(48, 5).ClearFormats()
$ws.Cells.Item(48, 7).NumberFormat = "@"
$ws.Cells.Item(48, 7).Value = "6"
$ws.Cells.Item(48, 7).ClearFormats()

# Row 49
$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = "148.02%"
$ws.Cells.Item(49, 5).ClearFormats()
$ws.Cells.Item(49, 7).NumberFormat = "@"
$ws.Cells.Item(49, 7).Value = "6"
$ws.Cells.Item(49, 7).ClearFormats()

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.00002103"
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = "0.09%"
$ws.Cells.Item(50, 5).ClearFormats()
$ws.Cells.Item(50, 7).NumberFormat = "@"
$ws.Cells.Item(50, 7).Value = "6"
$ws.Cells.Item(50, 7).ClearFormats()

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.0002003"
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = "0.09%"
$ws.Cells.Item(51, 5).ClearFormats()
$ws.Cells.Item(51, 7).NumberFormat = "@"
$ws.Cells.Item(51, 7).Value = "6"
$ws.Cells.Item(51, 7).ClearFormats()
